$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$cellData = @(
    @{ Row = 1; Col = 1; Lines = @("12 x 68", "  6    8", "  ----", "1|    |", "2|    |") }
    @{ Row = 1; Col = 2; Lines = @("83 x 67", "  6    7", "  ----", "8|    |", "3|    |") }
    @{ Row = 1; Col = 3; Lines = @("41 x 58", "  5    8", "  ----", "4|    |", "1|    |") }
    @{ Row = 2; Col = 1; Lines = @("22 x 64", "  6    4", "  ----", "2|    |", "2|    |") }
    @{ Row = 2; Col = 2; Lines = @("44 x 42", "  4    2", "  ----", "4|    |", "4|    |") }
    @{ Row = 2; Col = 3; Lines = @("66 x 26", "  2    6", "  ----", "6|    |", "6|    |") }
    @{ Row = 3; Col = 1; Lines = @("75 x 88", "  8    8", "  ----", "7|    |", "5|    |") }
    @{ Row = 3; Col = 2; Lines = @("39 x 99", "  9    9", "  ----", "3|    |", "9|    |") }
    @{ Row = 3; Col = 3; Lines = @("36 x 89", "  8    9", "  ----", "3|    |", "6|    |") }
    @{ Row = 4; Col = 1; Lines = @("97 x 94", "  9    4", "  ----", "9|    |", "7|    |") }
    @{ Row = 4; Col = 2; Lines = @("33 x 50", "  5    0", "  ----", "3|    |", "3|    |") }
    @{ Row = 4; Col = 3; Lines = @("16 x 19", "  1    9", "  ----", "1|    |", "6|    |") }
    @{ Row = 5; Col = 1; Lines = @("37 x 87", "  8    7", "  ----", "3|    |", "7|    |") }
    @{ Row = 5; Col = 2; Lines = @("66 x 75", "  7    5", "  ----", "6|    |", "6|    |") }
    @{ Row = 5; Col = 3; Lines = @("44 x 29", "  2    9", "  ----", "4|    |", "4|    |") }
)

foreach ($item in $cellData) {
    $cell = $t.Cell($item.Row, $item.Col)
    $newText = [string]::Join($nl, $item.Lines)
    $cell.Range.Text = $newText
}

Write-Host "Updated $($cellData.Count) cells"
